# Auto-generated edit script: update cryptocurrency Price (D) and Volume(1h) (E) columns
# for rows 2-51, matching the 'Updated cryptos list' GitHub Actions commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.952.82"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.63%  '
$ws.Range("D3").Value = "'1.846.00"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.81%  '
$ws.Range("D4").Value = "'1.006"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.28%  '
$ws.Range("D5").Value = "'309.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.17%  '
$ws.Range("D6").Value = "'1.006"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.30%  '
$ws.Range("D7").Value = "'0.4687"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.45%  '
$ws.Range("D8").Value = "'0.3662"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.85%  '
$ws.Range("D9").Value = "'0.07155"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.83%  '
$ws.Range("D10").Value = "'0.9268"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.72%  '
$ws.Range("D11").Value = "'19.60"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.25%  '
$ws.Range("D12").Value = "'0.07704"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.09%  '
$ws.Range("D13").Value = "'1.858.37"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.39%  '
$ws.Range("D14").Value = "'5.285"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.22%  '
$ws.Range("E15").Value = '  +1.25%  '
$ws.Range("D16").Value = "'88.48"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.82%  '
$ws.Range("E17").Value = '  +0.44%  '
$ws.Range("D18").Value = "'0.000008623"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.45%  '
$ws.Range("D19").Value = "'1.007"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.49%  '
$ws.Range("D20").Value = "'26.974.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.52%  '
$ws.Range("D21").Value = "'14.44"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.58%  '
$ws.Range("D22").Value = "'5.023"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.07%  '
$ws.Range("E23").Value = '  +0.78%  '
$ws.Range("D24").Value = "'1.917"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.12%  '
$ws.Range("D25").Value = "'152.28"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.97%  '
$ws.Range("D26").Value = "'18.22"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.30%  '
$ws.Range("D27").Value = "'2.020"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.51%  '
$ws.Range("D28").Value = "'114.36"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.60%  '
$ws.Range("D29").Value = "'4.881"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.54%  '
$ws.Range("D30").Value = "'0.08858"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.62%  '
$ws.Range("D31").Value = "'3.219"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.96%  '
$ws.Range("D32").Value = "'1.177"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +6.05%  '
$ws.Range("D33").Value = "'0.7465"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.19%  '
$ws.Range("D34").Value = "'2.795"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.34%  '
$ws.Range("D35").Value = "'4.473"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.75%  '
$ws.Range("D36").Value = "'1.085"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.21%  '
$ws.Range("D37").Value = "'0.01944"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.69%  '
$ws.Range("D38").Value = "'2.954"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.62%  '
$ws.Range("D39").Value = "'0.05193"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.77%  '
$ws.Range("E40").Value = '  +2.11%  '
$ws.Range("D41").Value = "'6.912"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.27%  '
$ws.Range("D42").Value = "'0.1521"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'8.134"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.07%  '
$ws.Range("D44").Value = "'10.50"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.90%  '
$ws.Range("D45").Value = "'0.4698"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.40%  '
$ws.Range("D46").Value = "'1.007"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.49%  '
$ws.Range("D47").Value = "'100.25"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.62%  '
$ws.Range("D48").Value = "'1.607"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.77%  '
$ws.Range("D49").Value = "'65.58"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.74%  '
$ws.Range("D50").Value = "'0.06043"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.99%  '
$ws.Range("D51").Value = "'0.8913"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.10%  '
